$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D2").Value = 0.1613222894898886
$ws.Range("E2").Value = 0.03454938278284612
$ws.Range("F2").Value = 0.03538930310934287
$ws.Range("G2").Value = 0.008068640211402529
$ws.Range("K2").Value = 0.0689944770795466
$ws.Range("L2").Value = 0.07292487232605809
$ws.Range("N2").Value = 0.2918381091764473
$ws.Range("O2").Value = 0.0002482012054538961
$ws.Range("P2").Value = 0.01329872165042877
$ws.Range("Q2").Value = 0.07595247679809845
$ws.Range("S2").Value = 0.1867897678074022
$ws.Range("T2").Value = 0.008606011936145908
$ws.Range("U2").Value = 0.02012029251601316
$ws.Range("AB2").Value = 0.01386804625222533
$ws.Range("AG2").Value = 0.008029407658700421
$ws.Range("D3").Value = 0.05802179344130491
$ws.Range("F3").Value = 0.1027319545916575
$ws.Range("G3").Value = 0.06936186619066691
$ws.Range("H3").Value = 0.00813568316824179
$ws.Range("J3").Value = 0.01203717762725273
$ws.Range("K3").Value = 0.02052474806959182
$ws.Range("M3").Value = 0.06653965466167951
$ws.Range("N3").Value = 0.1309767560127608
$ws.Range("O3").Value = 0.1291422787560598
$ws.Range("P3").Value = 0.04136758143277507
$ws.Range("R3").Value = 0.164043052619832
$ws.Range("S3").Value = 0.001074340802754849
$ws.Range("T3").Value = 0.1414363406012869
$ws.Range("U3").Value = 0.01474517389367799
$ws.Range("V3").Value = 0.01425147006927756
$ws.Range("W3").Value = 0.002086990766585298
$ws.Range("Z3").Value = 0.004448927674873909
$ws.Range("AC3").Value = 0.01907420961972067
$ws.Range("D4").Value = 0.1374397035822815
$ws.Range("E4").Value = 0.01429807020385547
$ws.Range("F4").Value = 0.07129927378839013
$ws.Range("G4").Value = 0.05134227900106719
$ws.Range("H4").Value = 0.002024443417255639
$ws.Range("I4").Value = 0.004760666813119077
$ws.Range("J4").Value = 0.008892878864057508
$ws.Range("L4").Value = 0.01462144527752529
$ws.Range("M4").Value = 0.08859118314161028
$ws.Range("N4").Value = 0.1229201691174575
$ws.Range("O4").Value = 0.1032138269521188
$ws.Range("P4").Value = 0.05856599955561008
$ws.Range("R4").Value = 0.06833188123597166
$ws.Range("S4").Value = 0.006311812573164285
$ws.Range("T4").Value = 0.186982537983224
$ws.Range("V4").Value = 0.007584178363496657
$ws.Range("W4").Value = 0.003814444118624861
$ws.Range("AC4").Value = 0.0197951906585514
$ws.Range("AE4").Value = 0.01565666269226719
$ws.Range("AH4").Value = 0.01355335266035133
$ws.Range("D5").Value = 0.1715567336402594
$ws.Range("E5").Value = 0.01388751504788677
$ws.Range("F5").Value = 0.05654070755197792
$ws.Range("G5").Value = 0.019137581346545
$ws.Range("K5").Value = 0.06298853276693887
$ws.Range("L5").Value = 0.05011465824785134
$ws.Range("N5").Value = 0.2355300864120008
$ws.Range("P5").Value = 0.009076718324801347
$ws.Range("Q5").Value = 0.109910836525318
$ws.Range("R5").Value = 0.02894203456019351
$ws.Range("S5").Value = 0.1691530488526429
$ws.Range("T5").Value = 0.02022536155419635
$ws.Range("U5").Value = 0.03274413800809134
$ws.Range("Z5").Value = 0.008401766306039705
$ws.Range("AC5").Value = 0.01179028085525688
$ws.Range("D6").Value = 0.1617311078452908
$ws.Range("F6").Value = 0.0857965503641495
$ws.Range("G6").Value = 0.09142081877735227
$ws.Range("I6").Value = 0.01382405660982873
$ws.Range("J6").Value = 0.008264340659643229
$ws.Range("K6").Value = 0.003867195268047824
$ws.Range("M6").Value = 0.03841562329792623
$ws.Range("N6").Value = 0.1720854398611379
$ws.Range("O6").Value = 0.06697798134333698
$ws.Range("P6").Value = 0.03686668816274159
$ws.Range("R6").Value = 0.08788013089811365
$ws.Range("S6").Value = 0.01179446964685298
$ws.Range("T6").Value = 0.1467798036531358
$ws.Range("W6").Value = 0.01035440073820133
$ws.Range("Z6").Value = 0.01668491778196502
$ws.Range("AC6").Value = 0.02804792327041222
$ws.Range("AE6").Value = 0.01286236971288941
$ws.Range("AH6").Value = 0.006346182108974349

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D2").Value = 0.1613222894898886
$ws.Range("E2").Value = 0.1958716722727348
$ws.Range("F2").Value = 0.2312609753820776
$ws.Range("G2").Value = 0.2393296155934801
$ws.Range("H2").Value = 0.2393296155934801
$ws.Range("I2").Value = 0.2393296155934801
$ws.Range("J2").Value = 0.2393296155934801
$ws.Range("K2").Value = 0.3083240926730267
$ws.Range("L2").Value = 0.3812489649990848
$ws.Range("M2").Value = 0.3812489649990848
$ws.Range("N2").Value = 0.6730870741755322
$ws.Range("O2").Value = 0.673335275380986
$ws.Range("P2").Value = 0.6866339970314148
$ws.Range("Q2").Value = 0.7625864738295133
$ws.Range("R2").Value = 0.7625864738295133
$ws.Range("S2").Value = 0.9493762416369155
$ws.Range("T2").Value = 0.9579822535730613
$ws.Range("U2").Value = 0.9781025460890744
$ws.Range("V2").Value = 0.9781025460890744
$ws.Range("W2").Value = 0.9781025460890744
$ws.Range("X2").Value = 0.9781025460890744
$ws.Range("Y2").Value = 0.9781025460890744
$ws.Range("Z2").Value = 0.9781025460890744
$ws.Range("AA2").Value = 0.9781025460890744
$ws.Range("AB2").Value = 0.9919705923412998
$ws.Range("AC2").Value = 0.9919705923412998
$ws.Range("AD2").Value = 0.9919705923412998
$ws.Range("AE2").Value = 0.9919705923412998
$ws.Range("AF2").Value = 0.9919705923412998
$ws.Range("D3").Value = 0.05802179344130491
$ws.Range("E3").Value = 0.05802179344130491
$ws.Range("F3").Value = 0.1607537480329624
$ws.Range("G3").Value = 0.2301156142236294
$ws.Range("H3").Value = 0.2382512973918711
$ws.Range("I3").Value = 0.2382512973918711
$ws.Range("J3").Value = 0.2502884750191239
$ws.Range("K3").Value = 0.2708132230887157
$ws.Range("L3").Value = 0.2708132230887157
$ws.Range("M3").Value = 0.3373528777503952
$ws.Range("N3").Value = 0.468329633763156
$ws.Range("O3").Value = 0.5974719125192158
$ws.Range("P3").Value = 0.6388394939519909
$ws.Range("Q3").Value = 0.6388394939519909
$ws.Range("R3").Value = 0.8028825465718229
$ws.Range("S3").Value = 0.8039568873745777
$ws.Range("T3").Value = 0.9453932279758646
$ws.Range("U3").Value = 0.9601384018695426
$ws.Range("V3").Value = 0.9743898719388201
$ws.Range("W3").Value = 0.9764768627054053
$ws.Range("X3").Value = 0.9764768627054053
$ws.Range("Y3").Value = 0.9764768627054053
$ws.Range("Z3").Value = 0.9809257903802793
$ws.Range("AA3").Value = 0.9809257903802793
$ws.Range("AB3").Value = 0.9809257903802793
$ws.Range("AC3").Value = 0.9999999999999999
$ws.Range("AD3").Value = 0.9999999999999999
$ws.Range("AE3").Value = 0.9999999999999999
$ws.Range("AF3").Value = 0.9999999999999999
$ws.Range("AG3").Value = 0.9999999999999999
$ws.Range("AH3").Value = 0.9999999999999999
$ws.Range("AI3").Value = 0.9999999999999999
$ws.Range("D4").Value = 0.1374397035822815
$ws.Range("E4").Value = 0.151737773786137
$ws.Range("F4").Value = 0.2230370475745271
$ws.Range("G4").Value = 0.2743793265755943
$ws.Range("H4").Value = 0.2764037699928499
$ws.Range("I4").Value = 0.281164436805969
$ws.Range("J4").Value = 0.2900573156700265
$ws.Range("K4").Value = 0.2900573156700265
$ws.Range("L4").Value = 0.3046787609475518
$ws.Range("M4").Value = 0.3932699440891621
$ws.Range("N4").Value = 0.5161901132066196
$ws.Range("O4").Value = 0.6194039401587383
$ws.Range("P4").Value = 0.6779699397143484
$ws.Range("Q4").Value = 0.6779699397143484
$ws.Range("R4").Value = 0.74630182095032
$ws.Range("S4").Value = 0.7526136335234843
$ws.Range("T4").Value = 0.9395961715067083
$ws.Range("U4").Value = 0.9395961715067083
$ws.Range("V4").Value = 0.9471803498702049
$ws.Range("W4").Value = 0.9509947939888298
$ws.Range("X4").Value = 0.9509947939888298
$ws.Range("Y4").Value = 0.9509947939888298
$ws.Range("Z4").Value = 0.9509947939888298
$ws.Range("AA4").Value = 0.9509947939888298
$ws.Range("AB4").Value = 0.9509947939888298
$ws.Range("AC4").Value = 0.9707899846473812
$ws.Range("AD4").Value = 0.9707899846473812
$ws.Range("AE4").Value = 0.9864466473396484
$ws.Range("AF4").Value = 0.9864466473396484
$ws.Range("AG4").Value = 0.9864466473396484
$ws.Range("AH4").Value = 0.9999999999999997
$ws.Range("AI4").Value = 0.9999999999999997
$ws.Range("D5").Value = 0.1715567336402594
$ws.Range("E5").Value = 0.1854442486881462
$ws.Range("F5").Value = 0.2419849562401241
$ws.Range("G5").Value = 0.2611225375866691
$ws.Range("H5").Value = 0.2611225375866691
$ws.Range("I5").Value = 0.2611225375866691
$ws.Range("J5").Value = 0.2611225375866691
$ws.Range("K5").Value = 0.324111070353608
$ws.Range("L5").Value = 0.3742257286014593
$ws.Range("M5").Value = 0.3742257286014593
$ws.Range("N5").Value = 0.6097558150134601
$ws.Range("O5").Value = 0.6097558150134601
$ws.Range("P5").Value = 0.6188325333382615
$ws.Range("Q5").Value = 0.7287433698635795
$ws.Range("R5").Value = 0.757685404423773
$ws.Range("S5").Value = 0.9268384532764158
$ws.Range("T5").Value = 0.9470638148306121
$ws.Range("U5").Value = 0.9798079528387034
$ws.Range("V5").Value = 0.9798079528387034
$ws.Range("W5").Value = 0.9798079528387034
$ws.Range("X5").Value = 0.9798079528387034
$ws.Range("Y5").Value = 0.9798079528387034
$ws.Range("Z5").Value = 0.9882097191447431
$ws.Range("AA5").Value = 0.9882097191447431
$ws.Range("AB5").Value = 0.9882097191447431
$ws.Range("D6").Value = 0.1617311078452908
$ws.Range("E6").Value = 0.1617311078452908
$ws.Range("F6").Value = 0.2475276582094403
$ws.Range("G6").Value = 0.3389484769867926
$ws.Range("H6").Value = 0.3389484769867926
$ws.Range("I6").Value = 0.3527725335966213
$ws.Range("J6").Value = 0.3610368742562645
$ws.Range("K6").Value = 0.3649040695243123
$ws.Range("L6").Value = 0.3649040695243123
$ws.Range("M6").Value = 0.4033196928222386
$ws.Range("N6").Value = 0.5754051326833765
$ws.Range("O6").Value = 0.6423831140267134
$ws.Range("P6").Value = 0.6792498021894551
$ws.Range("Q6").Value = 0.6792498021894551
$ws.Range("R6").Value = 0.7671299330875687
$ws.Range("S6").Value = 0.7789244027344217
$ws.Range("T6").Value = 0.9257042063875575
$ws.Range("U6").Value = 0.9257042063875575
$ws.Range("V6").Value = 0.9257042063875575
$ws.Range("W6").Value = 0.9360586071257588
$ws.Range("X6").Value = 0.9360586071257588
$ws.Range("Y6").Value = 0.9360586071257588
$ws.Range("Z6").Value = 0.9527435249077238
$ws.Range("AA6").Value = 0.9527435249077238
$ws.Range("AB6").Value = 0.9527435249077238
$ws.Range("AC6").Value = 0.9807914481781361
$ws.Range("AD6").Value = 0.9807914481781361
$ws.Range("AE6").Value = 0.9936538178910256
$ws.Range("AF6").Value = 0.9936538178910256
$ws.Range("AG6").Value = 0.9936538178910256
$ws.Range("AH6").Value = 0.9999999999999999
$ws.Range("AI6").Value = 0.9999999999999999

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.6730870741755322
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.5974719125192158
$ws.Range("G3").Value = 13
$ws.Range("F4").Value = 0.5161901132066196
$ws.Range("F5").Value = 0.6097558150134601
$ws.Range("F6").Value = 0.5754051326833765

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value = 16
$ws.Range("F2").Value = 0.7625864738295133
$ws.Range("G2").Value = 15
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.8028825465718229
$ws.Range("G3").Value = 16
$ws.Range("F4").Value = 0.74630182095032
$ws.Range("F5").Value = 0.7287433698635795
$ws.Range("D6").Value = 17
$ws.Range("F6").Value = 0.7671299330875687
$ws.Range("G6").Value = 16

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F2").Value = 0.9493762416369155
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.8028825465718229
$ws.Range("G3").Value = 16
$ws.Range("F4").Value = 0.9395961715067083
$ws.Range("F5").Value = 0.9268384532764158
$ws.Range("D6").Value = 19
$ws.Range("F6").Value = 0.9257042063875575
$ws.Range("G6").Value = 18

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F2").Value = 0.9493762416369155
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.9453932279758646
$ws.Range("G3").Value = 18
$ws.Range("F4").Value = 0.9395961715067083
$ws.Range("F5").Value = 0.9268384532764158
$ws.Range("F6").Value = 0.9257042063875575

Write-Output "Applied 274 cell updates"
